$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (FarsightVisionCircle) becomes a "Limited"/deconstructible class that
# deconstructs itself: Has Deconstruct? = Yes, Lifetime = Limited,
# Deconstructed At = Self. Also pick up the highlighted ("Limited" row) formatting
# by copying it from an existing row that already has that look (row 15).
$ws.Range("A15:F15").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)

$ws.Range("D2").Value = "Yes"
$ws.Range("E2").Value = "Limited"
$ws.Range("F2").Value = "Self"

# Row 17 (LimitedTable) lifetime corrected from Zombie to Limited.
$ws.Range("E17").Value = "Limited"

# Restore the selection state to match the saved view.
$ws.Range("I7").Select()
